$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Insert a new column before the existing "value"/year column (current column B),
# pushing the old data into column C.
$ws.Columns("B").Insert()

# Header row: new column B is the "budget-type" header.
$ws.Range("B1").Value = "budget-type"

# Fill the new column B for every data row (2..112) with "budget".
$lastRow = $ws.Cells(1,1).End(4).Row
$ws.Range("B2:B" + $lastRow).Value = "budget"
